$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 7 de Mayo de 2020 a las 11:04
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 11:04"

# Row 32: Israel
$ws.Range("B32").Value = 16346
$ws.Range("C32").Value = 36
$ws.Range("D32").Value = 10737
$ws.Range("E32").Value = 5370
$ws.Range("F32").Value = 83

# Row 33: Austria
$ws.Range("A33").Value = "Austria"
$ws.Range("B33").Value = 15752
$ws.Range("C33").Value = 68
$ws.Range("D33").Value = 13698
$ws.Range("E33").Value = 1445
$ws.Range("F33").Value = 92
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 609

# Row 34: Emiratos Arabes Unidos
$ws.Range("A34").Value = "Emiratos Arabes Unidos"
$ws.Range("B34").Value = 15738
$ws.Range("D34").Value = 3359
$ws.Range("E34").Value = 12222
$ws.Range("F34").Value = 1
$ws.Range("H34").Value = 157

# Row 36: Polonia
$ws.Range("B36").Value = 14898
$ws.Range("C36").Value = 158
$ws.Range("E36").Value = 9299
$ws.Range("G36").Value = 4
$ws.Range("H36").Value = 737

# Row 40: Banglades
$ws.Range("B40").Value = 12425
$ws.Range("C40").Value = 706
$ws.Range("D40").Value = 1910
$ws.Range("E40").Value = 10329

# Row 53: Malasia
$ws.Range("B53").Value = 6467
$ws.Range("C53").Value = 39
$ws.Range("D53").Value = 4776
$ws.Range("E53").Value = 1584
$ws.Range("F53").Value = 19

# Row 69: Armenia
$ws.Range("B69").Value = 2884
$ws.Range("C69").Value = 102
$ws.Range("D69").Value = 1185
$ws.Range("E69").Value = 1657
$ws.Range("G69").Value = 2
$ws.Range("H69").Value = 42

# Row 149: Brunei
$ws.Range("B149").Value = 141
$ws.Range("C149").Value = 2
$ws.Range("E149").Value = 9

# Row 150: Benin
$ws.Range("A150").Value = "Benin"
$ws.Range("B150").Value = 140
$ws.Range("C150").Value = 44
$ws.Range("D150").Value = 53
$ws.Range("E150").Value = 85
$ws.Range("H150").Value = 2

# Row 151: Guayana Francesa
$ws.Range("A151").Value = "Guayana Francesa"
$ws.Range("B151").Value = 138
$ws.Range("D151").Value = 112
$ws.Range("E151").Value = 25
$ws.Range("H151").Value = 1

# Row 152: Togo
$ws.Range("A152").Value = "Togo"
$ws.Range("B152").Value = 128
$ws.Range("D152").Value = 77
$ws.Range("E152").Value = 42
$ws.Range("H152").Value = 9

# Row 153: Suazilandia
$ws.Range("A153").Value = "Suazilandia"
$ws.Range("B153").Value = 123
$ws.Range("D153").Value = 12
$ws.Range("E153").Value = 109
$ws.Range("F153").Value = 0
$ws.Range("H153").Value = 2

# Row 154: Camboya
$ws.Range("A154").Value = "Camboya"
$ws.Range("B154").Value = 122
$ws.Range("D154").Value = 120
$ws.Range("E154").Value = 2
$ws.Range("F154").Value = 1
$ws.Range("H154").Value = 0

# Row 155: Bermudas
$ws.Range("A155").Value = "Bermudas"
$ws.Range("B155").Value = 118
$ws.Range("D155").Value = 59
$ws.Range("E155").Value = 52
$ws.Range("F155").Value = 4
$ws.Range("H155").Value = 7

# Row 156: Trinidad yTobago
$ws.Range("A156").Value = "Trinidad yTobago"
$ws.Range("B156").Value = 116
$ws.Range("D156").Value = 103
$ws.Range("E156").Value = 5
$ws.Range("H156").Value = 8

# Row 157: Haiti
$ws.Range("A157").Value = "Haiti"
$ws.Range("D157").Value = 10
$ws.Range("E157").Value = 79
$ws.Range("F157").Value = 0
$ws.Range("H157").Value = 12

# Row 158: Aruba
$ws.Range("A158").Value = "Aruba"
$ws.Range("B158").Value = 101
$ws.Range("D158").Value = 89
$ws.Range("E158").Value = 10
$ws.Range("F158").Value = 4
$ws.Range("H158").Value = 2

# Row 159: Uganda
$ws.Range("A159").Value = "Uganda"
$ws.Range("B159").Value = 100
$ws.Range("D159").Value = 55
$ws.Range("E159").Value = 45

# Row 160: Nepal
$ws.Range("A160").Value = "Nepal"
$ws.Range("B160").Value = 99
$ws.Range("D160").Value = 22
$ws.Range("E160").Value = 77
$ws.Range("H160").Value = 0

# Row 191: Nueva Caledonia
$ws.Range("A191").Value = "Nueva Caledonia"
$ws.Range("D191").Value = 18
$ws.Range("H191").Value = 0

# Row 192: Belice
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

# Row 198: Dominica
$ws.Range("A198").Value = "Dominica"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 0

# Row 199: Curazao
$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 13
$ws.Range("H199").Value = 1
